$wb = $excel.ActiveWorkbook

# --- Hoja1: update the two rate lines inside the "Conversion del dia" note (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")

$check = [string]([char]0x2705)   # U+2705 check-mark emoji, built via code point so this file stays plain ASCII

$oldText = $cellA1.Value2
$oldFragment = $check + " 1000 Bs = 10.05 = 42995.58 pesos`n" + $check + " 42995.58 pesos = 9.99 = 927.56 Bs"
$newFragment = $check + " 1000 Bs = 9.59 = 41040.62 pesos`n" + $check + " 41040.62 pesos = 9.54 = 945.37 Bs"

if ($oldText.Contains($oldFragment)) {
    $newText = $oldText.Replace($oldFragment, $newFragment)
    $cellA1.Value = $newText
}

# --- tasas: refresh the automatically-pulled rate figures (columns N/O, rows 10 & 12) ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 104.265
$wsTasas.Range("O10").Value = 4279.1
$wsTasas.Range("N12").Value = 4299.99
$wsTasas.Range("O12").Value = 99.05
